# Rename "wt" -> "wt_log2_expression" and "dcin5" -> "dcin5_log2_expression",
# then make "dcin5_log2_expression" the active sheet (was "optimization_parameters").

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("wt").Name = "wt_log2_expression"
$wb.Worksheets.Item("dcin5").Name = "dcin5_log2_expression"

$wb.Worksheets.Item("dcin5_log2_expression").Activate()
